# Apply the edits described by the diff:
# - B3: literal value 100 -> 69.067664339705
# - B4: formula changes from B2/(12*150) to B2/(12*95)
# - B5: formula changes from B3/(12*150) to B3/(12*95)
# - B6: literal value 0.077492756979866 -> 0.082048508723404
# - B7: literal value 19.675137572927 -> 15.507313055867

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 literal value change
$ws.Range("B3").Value = 69.067664339705

# B4 and B5 formulas now divide by 95 instead of 150
$ws.Range("B4").Formula = "=B2/(12*95)"
$ws.Range("B5").Formula = "=B3/(12*95)"

# B6 and B7 are plain literal (non-formula) values in the sheet
$ws.Range("B6").Value = 0.082048508723404
$ws.Range("B7").Value = 15.507313055867

$wb.Save()
